# Function_Descriptions.xlsx -- "Conductivity & Wrapping up Paper 1"
#
# 1. Add two explanatory comments (B17 / D17) to the existing find_time row.
# 2. Insert a new table row (new row 35) documenting
#    get_energy_distribution_moment, and shift everything below down by one
#    (which is why every comment anchored at B35/D35 and below moves down
#    one row in the diff).
# 3. Turn C34 (get_time_avg_time_series_data) into a hyperlink.
# 4. A handful of cosmetic tweaks: a couple of row heights, the width of
#    column C, and the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New comments on B17 / D17 (find_time's Output / Input cells)
# ---------------------------------------------------------------------
$cB17 = $ws.Range("B17").AddComment()
$cB17.Text("Nithin Sivadas:" + [char]10 + "% Output:" + [char]10 + "%  timeNo: The index of the time array that points " + [char]10 + "%             to the time specified by thisTime ")

$cD17 = $ws.Range("D17").AddComment()
$cD17.Text("Nithin Sivadas:" + [char]10 + "% Input" + [char]10 + "%  time     : 1-D time Array " + [char]10 + "%  thisTime : String identifying the time to be found " + [char]10 + "%             from the array time '26 Mar 2008 11:00'")

# ---------------------------------------------------------------------
# 2. Insert the new row for get_energy_distribution_moment above the old
#    row 35 (isThereNAN / totalNAN), which pushes it (and everything
#    after it) down to row 36, 37, ...
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$ws.Rows(35).Insert()
$lo.Resize($ws.Range("A1:F53"))

$ws.Range("C35").Value2 = "get_energy_distribution_moment"
$ws.Range("D35").Value2 = "flux,energyBin"
$ws.Range("B35").Value2 = "mean, median"
$ws.Range("E35").Value2 = "Calculates the mean and median of an input energy distribution"
$ws.Range("F35").Value2 = 42927

$ws.Range("A35:F35").WrapText = $true
$ws.Range("A35:F35").VerticalAlignment = -4108
$ws.Range("A35:F35").Borders.LineStyle = 1

$ws.Range("A35").HorizontalAlignment = -4108
$ws.Range("B35").HorizontalAlignment = -4152
$ws.Range("D35").HorizontalAlignment = -4131
$ws.Range("F35").HorizontalAlignment = -4108
$ws.Range("F35").NumberFormat = "m/d/yyyy"

$ws.Range("B35,C35,D35,F35").Font.Color = 8355711

$ws.Rows(35).RowHeight = 29

# ---------------------------------------------------------------------
# 3. Hyperlink on C34 (get_time_avg_time_series_data)
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C34"), "https://github.com/nithin-sivadas/energy-height-conversion") | Out-Null

# ---------------------------------------------------------------------
# 4. Cosmetic tweaks
# ---------------------------------------------------------------------
$ws.Rows(19).RowHeight = 58
$ws.Rows(20).RowHeight = 29
$ws.Rows(22).RowHeight = 29
$ws.Rows(30).RowHeight = 72.5

$ws.Columns("C").ColumnWidth = 25.8

$ws.Range("E17").Select()
